$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 322.4019826666667
$ws.Range("H2").Value = 967.205948
$ws.Range("I2").Value = 0.795674507658366
$ws.Range("J2").Value = 0.7956745076583662
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 152.2721466666667
$ws.Range("N2").Value = 456.81644
$ws.Range("O2").Value = 0.6052823830012941
$ws.Range("P2").Value = 0.6052823830012942
$ws.Range("Q2").Value = 49092.84199024279
$ws.Range("R2").Value = 441835.5779121851
$ws.Range("S2").Value = 0.4816077620888372
$ws.Range("T2").Value = 0.4816077620888374
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 322.4019826666667
$ws.Range("H3").Value = 967.205948
$ws.Range("I3").Value = 0.795674507658366
$ws.Range("J3").Value = 0.7956745076583662
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 66.66193
$ws.Range("N3").Value = 199.98579
$ws.Range("O3").Value = 0.2649814344194714
$ws.Range("P3").Value = 0.2649814344194714
$ws.Range("Q3").Value = 21491.93840038655
$ws.Range("R3").Value = 193427.4456034789
$ws.Range("S3").Value = 0.2108389723703205
$ws.Range("T3").Value = 0.2108389723703206
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 322.4019826666667
$ws.Range("H4").Value = 967.205948
$ws.Range("I4").Value = 0.795674507658366
$ws.Range("J4").Value = 0.7956745076583662
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.638001
$ws.Range("N4").Value = 97.91400300000001
$ws.Range("O4").Value = 0.1297361825792344
$ws.Range("P4").Value = 0.1297361825792344
$ws.Range("Q4").Value = 10522.55623267665
$ws.Range("R4").Value = 94703.00609408985
$ws.Range("S4").Value = 0.1032277731992082
$ws.Range("T4").Value = 0.1032277731992082
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 57.69151966666666
$ws.Range("H5").Value = 173.074559
$ws.Range("I5").Value = 0.1423802394983967
$ws.Range("J5").Value = 0.1423802394983967
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 152.2721466666667
$ws.Range("N5").Value = 456.81644
$ws.Range("O5").Value = 0.6052823830012941
$ws.Range("P5").Value = 0.6052823830012942
$ws.Range("Q5").Value = 8784.811544105551
$ws.Range("R5").Value = 79063.30389694996
$ws.Range("S5").Value = 0.08618025065588453
$ws.Range("T5").Value = 0.08618025065588455
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 57.69151966666666
$ws.Range("H6").Value = 173.074559
$ws.Range("I6").Value = 0.1423802394983967
$ws.Range("J6").Value = 0.1423802394983967
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 66.66193
$ws.Range("N6").Value = 199.98579
$ws.Range("O6").Value = 0.2649814344194714
$ws.Range("P6").Value = 0.2649814344194714
$ws.Range("Q6").Value = 3845.828045612956
$ws.Range("R6").Value = 34612.45241051661
$ws.Range("S6").Value = 0.03772812009527304
$ws.Range("T6").Value = 0.03772812009527304
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 57.69151966666666
$ws.Range("H7").Value = 173.074559
$ws.Range("I7").Value = 0.1423802394983967
$ws.Range("J7").Value = 0.1423802394983967
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.638001
$ws.Range("N7").Value = 97.91400300000001
$ws.Range("O7").Value = 0.1297361825792344
$ws.Range("P7").Value = 0.1297361825792344
$ws.Range("Q7").Value = 1882.935876572186
$ws.Range("R7").Value = 16946.42288914968
$ws.Range("S7").Value = 0.01847186874723911
$ws.Range("T7").Value = 0.01847186874723911
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 25.09980166666667
$ws.Range("H8").Value = 75.29940500000001
$ws.Range("I8").Value = 0.06194525284323717
$ws.Range("J8").Value = 0.06194525284323719
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 152.2721466666667
$ws.Range("N8").Value = 456.81644
$ws.Range("O8").Value = 0.6052823830012941
$ws.Range("P8").Value = 0.6052823830012942
$ws.Range("Q8").Value = 3822.000680690911
$ws.Range("R8").Value = 34398.00612621821
$ws.Range("S8").Value = 0.03749437025657228
$ws.Range("T8").Value = 0.0374943702565723
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 25.09980166666667
$ws.Range("H9").Value = 75.29940500000001
$ws.Range("I9").Value = 0.06194525284323717
$ws.Range("J9").Value = 0.06194525284323719
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 66.66193
$ws.Range("N9").Value = 199.98579
$ws.Range("O9").Value = 0.2649814344194714
$ws.Range("P9").Value = 0.2649814344194714
$ws.Range("Q9").Value = 1673.201221717217
$ws.Range("R9").Value = 15058.81099545495
$ws.Range("S9").Value = 0.01641434195387783
$ws.Range("T9").Value = 0.01641434195387783
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.09980166666667
$ws.Range("H10").Value = 75.29940500000001
$ws.Range("I10").Value = 0.06194525284323717
$ws.Range("J10").Value = 0.06194525284323719
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.638001
$ws.Range("N10").Value = 97.91400300000001
$ws.Range("O10").Value = 0.1297361825792344
$ws.Range("P10").Value = 0.1297361825792344
$ws.Range("Q10").Value = 819.2073518964685
$ws.Range("R10").Value = 7372.866167068216
$ws.Range("S10").Value = 0.008036540632787055
$ws.Range("T10").Value = 0.008036540632787057
